$d = $word.ActiveDocument
$apos = [char]0x2019

function Find-ParagraphByText($doc, $pattern) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like $pattern) {
            return $p
        }
    }
    return $null
}

# Replaces the entire contents of a paragraph (keeping its paragraph
# mark) with the given inner-paragraph OOXML (w:pPr / w:r ... ), so the
# emitted run/xml:space structure exactly matches what we want instead
# of whatever run-splitting/merging a plain text assignment would do.
function Set-ParagraphInnerXml($doc, $para, $innerXml) {
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End - 1   # exclude the paragraph mark itself
    # Always clear (even an already-collapsed range) - this also drops
    # any stray empty <w:r/> that InsertParagraphBefore() leaves behind.
    $doc.Range($pStart, $pEnd).Text = ""
    $insPos = $doc.Range($pStart, $pStart)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insPos.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) Add a new "Author" paragraph "Sarah Sant'Ana" right before the
#    existing "Olivia Cox" author paragraph.
# ---------------------------------------------------------------------
$oliviaPara = Find-ParagraphByText $d "*Olivia*Cox*"
if ($oliviaPara -ne $null) {
    $oliviaPara.Range.InsertParagraphBefore()
    $newAuthorPara = Find-ParagraphByText $d "*Olivia*Cox*"
    $authorInner = '<w:p><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">Sarah</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Sant' + $apos + 'Ana</w:t></w:r></w:p>'
    # the blank paragraph inserted right before "Olivia Cox" is the one
    # we just created - locate it via its position (previous paragraph).
    $idx = 0
    $targetIdx = -1
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like "*Olivia*Cox*") {
            $targetIdx = $idx - 1
            break
        }
    }
    $blankPara = $d.Paragraphs($targetIdx)
    Set-ParagraphInnerXml $d $blankPara $authorInner
}

# ---------------------------------------------------------------------
# 2) Update the date from 2024-08-21 to 2024-09-03.
# ---------------------------------------------------------------------
$datePara = Find-ParagraphByText $d "*2024-08-21*"
if ($datePara -ne $null) {
    $dateInner = '<w:p><w:pPr><w:pStyle w:val="Date"/></w:pPr><w:r><w:t xml:space="preserve">2024-09-03</w:t></w:r></w:p>'
    Set-ParagraphInnerXml $d $datePara $dateInner
}

# ---------------------------------------------------------------------
# 3) Add a new "BodyText" paragraph right before the
#    "Algorithm-guided risk monitoring..." paragraph.
# ---------------------------------------------------------------------
$algoPara = Find-ParagraphByText $d "*Algorithm-guided risk monitoring*"
if ($algoPara -ne $null) {
    $algoPara.Range.InsertParagraphBefore()
    $idx = 0
    $targetIdx = -1
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like "*Algorithm-guided risk monitoring*") {
            $targetIdx = $idx - 1
            break
        }
    }
    $blankBodyPara = $d.Paragraphs($targetIdx)
    $bodyInner = '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">High risk of overdose after initial period of remission.</w:t></w:r></w:p>'
    Set-ParagraphInnerXml $d $blankBodyPara $bodyInner
}

# ---------------------------------------------------------------------
# 4) Remove "even " from "People can comply with even highly effortful".
# ---------------------------------------------------------------------
$peoplePara = Find-ParagraphByText $d "*People can comply with even highly effortful*"
if ($peoplePara -ne $null) {
    $peopleInner = '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">People can comply with highly effortful sensing methods (e.g., 4 x daily EMA) while using substances</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">(Wyant et al. 2023; Jones et al. 2019)</w:t></w:r><w:r><w:t xml:space="preserve">.</w:t></w:r></w:p>'
    Set-ParagraphInnerXml $d $peoplePara $peopleInner
}

Write-Output "done"
